# [TASK] Update Wireframes and Filters
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# --- Row 108: end time of "Filterarten bestimmen und konzipieren" corrected ---
$ws.Range("K108").Value = 0.66666666666666663

# --- Row 109: new entry "Domaenenmodell ueberarbeiten" / "Iteration anfertigen" ---
# Pick up the I/J/K number-formatting (formula number format + filled/bordered time cells)
# from the row right above so the new cells match the rest of the table.
$ws.Range("I108:K108").Copy()
$ws.Range("I109:K109").PasteSpecial(-4122)

$ws.Range("A109").Value = 23
$ws.Range("B109").Value = "Interface Design"
$ws.Range("C109").Value = "Domaenenmodell ueberarbeiten"
$ws.Range("D109").Value = "[TASK]"
$ws.Range("E109").Value = "Iteration anfertigen"
$ws.Range("F109").Value = 44379
$ws.Range("G109").Value = 44359
$ws.Range("I109").Formula = "=ROUNDUP(((SUM(K109-J109)*24*60/60)/0.25),0)*0.25"
$ws.Range("J109").Value = 0.70833333333333337
$ws.Range("K109").Value = 0.71875

# --- Row 110: new entry "MockUps" / "MockUps Rezept Filter" ---
# Same I/J/K formatting as above, plus the bold-ish black-font styling (A/B/C/E)
# used by the other "MockUps" rows (105/106).
$ws.Range("I108:K108").Copy()
$ws.Range("I110:K110").PasteSpecial(-4122)
$ws.Range("A105:C105").Copy()
$ws.Range("A110:C110").PasteSpecial(-4122)
$ws.Range("E105").Copy()
$ws.Range("E110").PasteSpecial(-4122)

$ws.Range("A110").Value = 22
$ws.Range("B110").Value = "Interface Design"
$ws.Range("C110").Value = "MockUps"
$ws.Range("D110").Value = "[FEATURE]"
$ws.Range("E110").Value = "MockUps Rezept Filter"
$ws.Range("F110").Value = 44379
$ws.Range("G110").Value = 44359
$ws.Range("I110").Formula = "=ROUNDUP(((SUM(K110-J110)*24*60/60)/0.25),0)*0.25"
$ws.Range("J110").Value = 0.71875
$ws.Range("K110").Value = 0.79166666666666663

# --- sheet view: selection moved one row further down after the insert ---
$ws.Range("M109").Select() | Out-Null
